$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A279").Value = "IMX-USD"
$ws.Range("A280").Value = "TAO-USD"
$ws.Range("A281").Value = "MNT-USD"
